# Adds "Mid Paper 1/2/3" columns (G:I) to both marks sheets, carrying over
# the existing Paper 1/2/3 marks as the "Mid" term marks, and switches the
# active tab/selection over to the "Senior Five" sheet.

$wb = $excel.ActiveWorkbook

$wsSix  = $wb.Worksheets.Item("Senior Six")
$wsFive = $wb.Worksheets.Item("Senior Five")

function Add-MidPaperColumns($ws) {
    # Header row: G1/H1/I1 get the three new "Mid Paper n" headers, mirroring
    # the existing Paper 1/2/3 columns (D1/E1/F1) one column group over.
    $ws.Range("G1").Value = "Mid Paper 1"
    $ws.Range("H1").Value = "Mid Paper 2"
    $ws.Range("I1").Value = "Mid Paper 3"

    for ($r = 2; $r -le 6; $r++) {
        $p1 = $ws.Cells.Item($r, 4).Value2
        $p2 = $ws.Cells.Item($r, 5).Value2
        $p3 = $ws.Cells.Item($r, 6).Value2

        # Mirror D/E formatting onto G/H even when the source marks are
        # still blank, so the new "Mid Paper" cells exist on every row.
        $ws.Cells.Item($r, 7).Style = $ws.Cells.Item($r, 4).Style
        $ws.Cells.Item($r, 8).Style = $ws.Cells.Item($r, 5).Style

        if ($p1 -ne $null) {
            $ws.Cells.Item($r, 7).Value = $p1
        }
        if ($p2 -ne $null) {
            $ws.Cells.Item($r, 8).Value = $p2
        }
        if ($p3 -ne $null) {
            $ws.Cells.Item($r, 9).Value = $p3
        }
    }
}

Add-MidPaperColumns $wsSix
Add-MidPaperColumns $wsFive

# Senior Five's "Mid Paper 3" total came from the Senior Six sheet instead
# of its own Paper 3 column in the source edit - replicate that quirk.
$wsFive.Range("I3").Value = 82
$wsFive.Range("I4").Value = 76

# Move the active selection / tab over to "Senior Five".
$wsSix.Activate()
$wsSix.Range("G1").Select()

$wsFive.Activate()
$wsFive.Range("G11").Select()
